$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 7136248000
$ws.Range("C2").Value = "9811 Katy Fwy #100, Houston, TX 77024"

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 9056626552
$ws.Range("C3").Value = "324 Leaside Avenue Stoney Creek, Ontario L8E 2N7 Canada"

$ws.Range("A4").Value = 6
$ws.Range("B4").Value = 4178334565
$ws.Range("C4").Value = "3715 E Farm Road 94 Springfield, MO 65803"

$ws.Range("A5").Value = 32
$ws.Range("B5").Value = 8002226283
$ws.Range("C5").Value = "223 15th St NE Sioux Center, IA 51250 United States"

$ws.Range("A6").Value = 59
$ws.Range("B6").Value = 9725476020
$ws.Range("C6").Value = "725 E University Dr, McKinney, TX 75069"

$ws.Columns.Item(3).ColumnWidth = 49

$ws.Range("C11").Select()
